$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, pushing the existing data (old rows 34-96)
# down to rows 35-97.
$ws.Rows("34:34").Insert()

# Populate the newly inserted row 34 with the new price-report entry.
$ws.Cells.Item(34, 1).Value = 10
$ws.Cells.Item(34, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(34, 3).Value = "La Araucanía"
$ws.Cells.Item(34, 4).Value = 45014
$ws.Cells.Item(34, 5).Value = 9
$ws.Cells.Item(34, 6).Value = "Fruta"
$ws.Cells.Item(34, 7).Value = 100108
$ws.Cells.Item(34, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(34, 9).Value = 100108004
$ws.Cells.Item(34, 10).Value = "Papaya"
$ws.Cells.Item(34, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(34, 12).Value = "Primera"
$ws.Cells.Item(34, 13).Value = 110
$ws.Cells.Item(34, 14).Value = 28000
$ws.Cells.Item(34, 15).Value = 28000
$ws.Cells.Item(34, 16).Value = 28000
$ws.Cells.Item(34, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(34, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(34, 19).Value = 2800
$ws.Cells.Item(34, 20).Value = 10
